$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2131189209.519208
$ws.Range("E3").Value = 1996070030.511242
$ws.Range("E4").Value = 1225934278.599479
$ws.Range("E5").Value = 772476100.9553949
$ws.Range("E6").Value = 880716795.9077289
$ws.Range("E7").Value = 1294451877.349769
$ws.Range("E8").Value = 3273209521.527452
$ws.Range("E9").Value = 2761908704.267637
$ws.Range("E10").Value = 3437881270.23809
$ws.Range("E11").Value = 1824132549.423184
$ws.Range("E12").Value = 1342886505.16142
$ws.Range("E13").Value = 2525835536.001273
$ws.Range("E14").Value = 2544553617.619198
$ws.Range("E15").Value = 1426551350.42554
$ws.Range("E16").Value = 1252422602.001582
$ws.Range("E17").Value = 1580101624.698467
$ws.Range("E18").Value = 2493007861.836544
$ws.Range("E19").Value = 2336955899.780849
$ws.Range("E20").Value = 746618758.7015805
$ws.Range("E21").Value = 636156600.9461595
$ws.Range("E22").Value = 810928572.5470011
$ws.Range("E23").Value = 1896832779.482464
$ws.Range("E24").Value = 2283733516.651238
$ws.Range("E25").Value = 1494514463.522902
$ws.Range("E26").Value = 2264374430.546672
$ws.Range("E27").Value = 1626829036.042399
$ws.Range("E28").Value = 1921284947.617015
$ws.Range("E29").Value = 1568172438.185875
$ws.Range("E30").Value = 1084448908.864102
$ws.Range("E31").Value = 1677019772.88577
$ws.Range("E32").Value = 1969026529.05234
$ws.Range("E33").Value = 2122624846.609889
$ws.Range("E34").Value = 1719621813.176892
$ws.Range("E35").Value = 1812071094.593312
